# results with fixed workflow
#
# Both data sheets (NBR, BAR) get their "Cutoff"/"Reaction_number" table
# refreshed with new simulation results: the window of cutoffs now starts
# at 5 instead of 1 (column B), the reaction counts (column C) are
# recomputed, and the table shrinks from 19 data rows (A1:C20) down to
# 15 data rows (A1:C16) by dropping the last four rows.

$wb = $excel.ActiveWorkbook

$sheetData = @{
    "NBR" = @{
        B = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
        C = @(882,874,874,873,863,877,874,872,868,847,833,827,826,827,820)
    }
    "BAR" = @{
        B = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
        C = @(841,840,839,837,838,797,796,797,791,794,793,796,792,783,780)
    }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $sheetData.ContainsKey($name)) {
        continue
    }
    $bVals = $sheetData[$name].B
    $cVals = $sheetData[$name].C

    # Rows 2..16 (data rows 0..14) hold the refreshed values.
    for ($i = 0; $i -lt $bVals.Count; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $bVals[$i]
        $ws.Cells.Item($row, 3).Value = $cVals[$i]
    }

    # Drop the trailing four rows (old rows 17:20) so the table ends at row 16.
    $ws.Rows("17:20").Delete()
}
